$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 1.81
$ws.Range("K2").Value = 5.1
$ws.Range("N2").Value = 8
$ws.Range("S2").Value = 1.96
$ws.Range("U2").Value = 2.96
$ws.Range("W2").Value = 2.22
$ws.Range("AK2").Value = 1000
$ws.Range("F3").Value = 2.86
$ws.Range("G3").Value = 3.55
$ws.Range("H3").Value = 2.72
$ws.Range("I3").Value = 3.1
$ws.Range("J3").Value = 2.52
$ws.Range("K3").Value = 3.05
$ws.Range("M3").Value = 1.17
$ws.Range("N3").Value = 2.16
$ws.Range("O3").Value = 1.7
$ws.Range("P3").Value = 1.38
$ws.Range("Q3").Value = 2.8
$ws.Range("R3").Value = 1.13
$ws.Range("S3").Value = 1.02
$ws.Range("T3").Value = 2.32
$ws.Range("U3").Value = 1.61
$ws.Range("V3").Value = 1.47
$ws.Range("W3").Value = 1.39
$ws.Range("X3").Value = 8
$ws.Range("Y3").Value = 8.800000000000001
$ws.Range("Z3").Value = 17.5
$ws.Range("AA3").Value = 60
$ws.Range("AB3").Value = 9.6
$ws.Range("AC3").Value = 8.4
$ws.Range("AD3").Value = 18
$ws.Range("AE3").Value = 55
$ws.Range("AF3").Value = 21
$ws.Range("AG3").Value = 20
$ws.Range("AH3").Value = 30
$ws.Range("AI3").Value = 120
$ws.Range("AJ3").Value = 75
$ws.Range("AK3").Value = 65
$ws.Range("AL3").Value = 130
$ws.Range("AM3").Value = 350
$ws.Range("AN3").Value = 120
$ws.Range("AO3").Value = 90
$ws.Range("G4").Value = 6
$ws.Range("I4").Value = 1.94
$ws.Range("J4").Value = 3.35
$ws.Range("M4").Value = 1.09
$ws.Range("Q4").Value = 2.14
$ws.Range("S4").Value = 4
$ws.Range("V4").Value = 2.06
$ws.Range("P5").Value = 2.64
$ws.Range("Q5").Value = 1.56
$ws.Range("T5").Value = 1.95
$ws.Range("AJ5").Value = 400
$ws.Range("AK5").Value = 170
$ws.Range("AL5").Value = 130
$ws.Range("P6").Value = 1.68
$ws.Range("S7").Value = 2.48
$ws.Range("AI7").Value = 410
$ws.Range("AM7").Value = 430
$ws.Range("I10").Value = 4.3
$ws.Range("K10").Value = 3.15
